$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 57000
$ws.Range("J3").Value = 57000
$ws.Range("L3").Value = 57000
$ws.Range("N3").Value = -57228
$ws.Range("H17").Value = 2688.375
$ws.Range("J17").Value = 2688.375
$ws.Range("L17").Value = 8065.125
$ws.Range("N17").Value = -8401.125
$ws.Range("H81").Value = 20000
$ws.Range("I81").Value = 20000
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 20000
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -19002
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 20000
$ws.Range("I84").Value = 20000
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 60000
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -55008
$ws.Range("N84").ClearContents()
$ws.Range("H97").Value = 2354.8333
$ws.Range("I97").Value = 1000
$ws.Range("J97").Value = 2625.8
$ws.Range("K97").Value = 3000
$ws.Range("L97").Value = 7877.400000000001
$ws.Range("M97").Value = -2504
$ws.Range("N97").Value = -8869.400000000001
$ws.Range("H102").Value = 57000
$ws.Range("J102").Value = 57000
$ws.Range("L102").Value = 57000
$ws.Range("N102").Value = -63490
$ws.Range("H121").Value = 4766.077
$ws.Range("J121").Value = 4766.077
$ws.Range("L121").Value = 14298.231
$ws.Range("N121").Value = -17792.231
$ws.Range("H131").Value = 2399.5625
$ws.Range("I131").Value = 1324.6666
$ws.Range("J131").Value = 5624.25
$ws.Range("K131").Value = 3973.9998
$ws.Range("L131").Value = 16872.75
$ws.Range("M131").Value = 1066.0002
$ws.Range("N131").Value = -26952.75
$ws.Range("H135").Value = 2463.5334
$ws.Range("I135").Value = 748.5909
$ws.Range("K135").Value = 6737.3181
$ws.Range("M135").Value = -4202.3181
$ws.Range("H138").Value = 1962.58
$ws.Range("I138").Value = 1218.96
$ws.Range("J138").Value = 2210.4534
$ws.Range("K138").Value = 3656.88
$ws.Range("L138").Value = 6631.360199999999
$ws.Range("M138").Value = 1483.12
$ws.Range("N138").Value = -16911.3602
$ws.Range("H139").Value = 122895
$ws.Range("J139").Value = 122895
$ws.Range("L139").Value = 122895
$ws.Range("N139").Value = -133175

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 761736.4
$ws.Range("I2").Value = 875611.8
$ws.Range("K2").Value = 875611.8
$ws.Range("M2").Value = -875498.8
$ws.Range("H46").Value = 19142.334
$ws.Range("J46").Value = 19766.785
$ws.Range("L46").Value = 19766.785
$ws.Range("N46").Value = -20404.785
$ws.Range("H74").Value = 1335.7949
$ws.Range("I74").Value = 1152.1034
$ws.Range("J74").Value = 1868.5
$ws.Range("K74").Value = 1152.1034
$ws.Range("L74").Value = 1868.5
$ws.Range("M74").Value = -278.1034
$ws.Range("N74").Value = -3616.5
$ws.Range("H77").Value = 1335.7949
$ws.Range("I77").Value = 1152.1034
$ws.Range("J77").Value = 1868.5
$ws.Range("K77").Value = 5760.517
$ws.Range("L77").Value = 9342.5
$ws.Range("M77").Value = -1392.517
$ws.Range("N77").Value = -18078.5
$ws.Range("H102").Value = 654095.7
$ws.Range("I102").Value = 807060
$ws.Range("J102").Value = 3997.5
$ws.Range("K102").Value = 807060
$ws.Range("L102").Value = 3997.5
$ws.Range("M102").Value = -805438
$ws.Range("N102").Value = -7241.5
$ws.Range("H116").Value = 761736.4
$ws.Range("I116").Value = 875611.8
$ws.Range("K116").Value = 875611.8
$ws.Range("M116").Value = -873317.8
$ws.Range("H132").Value = 12916.579
$ws.Range("I132").Value = 19938.967
$ws.Range("J132").Value = 4543.731
$ws.Range("K132").Value = 59816.901
$ws.Range("L132").Value = 13631.193
$ws.Range("M132").Value = -57286.901
$ws.Range("N132").Value = -18691.193

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 761736.4
$ws.Range("I3").Value = 875611.8
$ws.Range("K3").Value = 875611.8
$ws.Range("M3").Value = -875497.8
$ws.Range("H86").Value = 4344.6665
$ws.Range("J86").Value = 5074.2856
$ws.Range("L86").Value = 5074.2856
$ws.Range("N86").Value = -7320.2856
$ws.Range("H89").Value = 4344.6665
$ws.Range("J89").Value = 5074.2856
$ws.Range("L89").Value = 25371.428
$ws.Range("N89").Value = -36603.428
$ws.Range("H105").Value = 2285.4
$ws.Range("I105").Value = 2175.5625
$ws.Range("K105").Value = 2175.5625
$ws.Range("M105").Value = -428.5625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2726.5
$ws.Range("I16").Value = 2643.75
$ws.Range("K16").Value = 2643.75
$ws.Range("M16").Value = -2356.75
$ws.Range("H113").Value = 2726.5
$ws.Range("I113").Value = 2643.75
$ws.Range("K113").Value = 2643.75
$ws.Range("M113").Value = -473.75
$ws.Range("H122").Value = 5027.8237
$ws.Range("I122").Value = 3065.5715
$ws.Range("J122").Value = 6401.4
$ws.Range("K122").Value = 9196.7145
$ws.Range("L122").Value = 19204.2
$ws.Range("M122").Value = -6746.7145
$ws.Range("N122").Value = -24104.2
$ws.Range("H141").Value = 83310.055
$ws.Range("J141").Value = 84663.82
$ws.Range("L141").Value = 84663.82
$ws.Range("N141").Value = -95023.82

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1205.2
$ws.Range("I22").Value = 250
$ws.Range("J22").Value = 1444
$ws.Range("K22").Value = 750
$ws.Range("L22").Value = 4332
$ws.Range("M22").Value = -581
$ws.Range("N22").Value = -4670
$ws.Range("H27").Value = 1205.2
$ws.Range("I27").Value = 250
$ws.Range("J27").Value = 1444
$ws.Range("K27").Value = 750
$ws.Range("L27").Value = 4332
$ws.Range("M27").Value = -648
$ws.Range("N27").Value = -4536
$ws.Range("H60").Value = 324
$ws.Range("I60").Value = 236
$ws.Range("J60").Value = 500
$ws.Range("K60").Value = 708
$ws.Range("L60").Value = 1500
$ws.Range("M60").Value = -457
$ws.Range("N60").Value = -2002
$ws.Range("H103").Value = 1750
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 1750
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 5250
$ws.Range("M103").ClearContents()
$ws.Range("N103").Value = -7008

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 20842778
$ws.Range("I102").Value = 35724544
$ws.Range("K102").Value = 35724544
$ws.Range("M102").Value = -35722922
$ws.Range("H113").Value = 3799.8333
$ws.Range("I113").Value = 3559.8
$ws.Range("K113").Value = 3559.8
$ws.Range("M113").Value = -1389.8
$ws.Range("H132").Value = 367762.28
$ws.Range("I132").Value = 127714.875
$ws.Range("J132").Value = 593689.25
$ws.Range("K132").Value = 383144.625
$ws.Range("L132").Value = 1781067.75
$ws.Range("M132").Value = -380614.625
$ws.Range("N132").Value = -1786127.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4891.087
$ws.Range("I7").Value = 2824.8333
$ws.Range("J7").Value = 7145.1816
$ws.Range("K7").Value = 2824.8333
$ws.Range("L7").Value = 7145.1816
$ws.Range("M7").Value = -2712.8333
$ws.Range("N7").Value = -7369.1816
$ws.Range("H22").Value = 1154.6364
$ws.Range("J22").Value = 1333.3334
$ws.Range("L22").Value = 1333.3334
$ws.Range("N22").Value = -1923.3334
$ws.Range("H27").Value = 1154.6364
$ws.Range("J27").Value = 1333.3334
$ws.Range("L27").Value = 1333.3334
$ws.Range("N27").Value = -1547.3334
$ws.Range("H40").Value = 3088.4614
$ws.Range("I40").Value = 3088.4614
$ws.Range("K40").Value = 3088.4614
$ws.Range("M40").Value = -2952.4614
$ws.Range("H82").Value = 2842399.5
$ws.Range("I82").Value = 6250619
$ws.Range("K82").Value = 6250619
$ws.Range("M82").Value = -6250258
$ws.Range("H85").Value = 2842399.5
$ws.Range("I85").Value = 6250619
$ws.Range("K85").Value = 6250619
$ws.Range("M85").Value = -6249371
$ws.Range("H122").Value = 9480.4
$ws.Range("I122").Value = 3945.8
$ws.Range("K122").Value = 11837.4
$ws.Range("M122").Value = -9387.400000000001
$ws.Range("H126").Value = 4891.087
$ws.Range("I126").Value = 2824.8333
$ws.Range("J126").Value = 7145.1816
$ws.Range("K126").Value = 8474.499899999999
$ws.Range("L126").Value = 21435.5448
$ws.Range("M126").Value = -6004.499899999999
$ws.Range("N126").Value = -26375.5448
$ws.Range("H132").Value = 3750.027
$ws.Range("J132").Value = 3833.6
$ws.Range("L132").Value = 11500.8
$ws.Range("N132").Value = -16560.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 79085.78
$ws.Range("J80").Value = 83475.5
$ws.Range("L80").Value = 83475.5
$ws.Range("N80").Value = -85471.5
$ws.Range("H83").Value = 79085.78
$ws.Range("J83").Value = 83475.5
$ws.Range("L83").Value = 250426.5
$ws.Range("N83").Value = -260410.5
$ws.Range("H122").Value = 4734.5625
$ws.Range("I122").Value = 4734.5625
$ws.Range("K122").Value = 14203.6875
$ws.Range("M122").Value = -11753.6875
$ws.Range("H132").Value = 23814100
$ws.Range("I132").Value = 1413.1333
$ws.Range("K132").Value = 4239.3999
$ws.Range("M132").Value = -1709.3999
